# ---- Auto-generated PowerShell COM-interop edit script ----
# Adds a "2022-Q1" worksheet (13 fund rows) between "2020-Q4" and "总计",
# and inserts a matching 2022-Q1 summary row at the top of "总计".
$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2020-Q4")

# 1. Create the new sheet by copying "总计" (so the header / index-column cells
#    inherit the same style, s=2) and drop it right after "2020-Q4".
$totalSheet.Copy($null, $q4Sheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"
$newSheet.Cells.Clear()

# Re-apply the header style (copied from 总计!B1) across B1:H1 ...
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
# ... and the index-column style (copied from 总计!A2) down A2:A14.
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

# 2. Header row.
$newSheet.Range("B1").Value = '基金代码'
$newSheet.Range("C1").Value = '基金名称'
$newSheet.Range("D1").Value = '基金规模'
$newSheet.Range("E1").Value = '股票总仓位'
$newSheet.Range("F1").Value = '仓位占比'
$newSheet.Range("G1").Value = '持有市值(亿元)'
$newSheet.Range("H1").Value = '仓位排名'

# 3. Data rows 2..14: one literal assignment per cell so every value keeps its
#    exact text (fund codes such as 003318 must not be coerced into numbers).
$newSheet.Range("B2:G14").NumberFormat = "@"
$newSheet.Range("H2:H14").NumberFormat = "General"

# row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = '003318'
$newSheet.Range("C2").Value = '景顺长城中证500行业中性低波动指数'
$newSheet.Range("D2").Value = '13.99'
$newSheet.Range("E2").Value = '93.88'
$newSheet.Range("F2").Value = '2.48'
$newSheet.Range("G2").Value = '0.3470'
$newSheet.Range("H2").Value = 1

# row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = '008114'
$newSheet.Range("C3").Value = '天弘中证红利低波动100指数A'
$newSheet.Range("D3").Value = '3.16'
$newSheet.Range("E3").Value = '92.60'
$newSheet.Range("F3").Value = '2.40'
$newSheet.Range("G3").Value = '0.0758'
$newSheet.Range("H3").Value = 2

# row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = '008115'
$newSheet.Range("C4").Value = '天弘中证红利低波动100指数C'
$newSheet.Range("D4").Value = '2.37'
$newSheet.Range("E4").Value = '92.60'
$newSheet.Range("F4").Value = '2.40'
$newSheet.Range("G4").Value = '0.0569'
$newSheet.Range("H4").Value = 2

# row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = '000916'
$newSheet.Range("C5").Value = '前海开源股息率100强等权重股票'
$newSheet.Range("D5").Value = '3.07'
$newSheet.Range("E5").Value = '90.96'
$newSheet.Range("F5").Value = '1.76'
$newSheet.Range("G5").Value = '0.0540'
$newSheet.Range("H5").Value = 2

# row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = '519619'
$newSheet.Range("C6").Value = '银河君荣灵活配置混合A'
$newSheet.Range("D6").Value = '0.80'
$newSheet.Range("E6").Value = '93.36'
$newSheet.Range("F6").Value = '6.00'
$newSheet.Range("G6").Value = '0.0480'
$newSheet.Range("H6").Value = 5

# row 7
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = '519621'
$newSheet.Range("C7").Value = '银河君荣灵活配置混合I'
$newSheet.Range("D7").Value = '0.80'
$newSheet.Range("E7").Value = '93.36'
$newSheet.Range("F7").Value = '6.00'
$newSheet.Range("G7").Value = '0.0480'
$newSheet.Range("H7").Value = 5

# row 8
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = '515100'
$newSheet.Range("C8").Value = '景顺长城中证红利低波动100ETF'
$newSheet.Range("D8").Value = '1.25'
$newSheet.Range("E8").Value = '97.96'
$newSheet.Range("F8").Value = '2.59'
$newSheet.Range("G8").Value = '0.0324'
$newSheet.Range("H8").Value = 2

# row 9
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = '512260'
$newSheet.Range("C9").Value = '华安中证500行业中性低波动ETF'
$newSheet.Range("D9").Value = '1.17'
$newSheet.Range("E9").Value = '96.94'
$newSheet.Range("F9").Value = '2.56'
$newSheet.Range("G9").Value = '0.0300'
$newSheet.Range("H9").Value = 1

# row 10
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = '159990'
$newSheet.Range("C10").Value = '银华巨潮小盘价值ETF'
$newSheet.Range("D10").Value = '1.06'
$newSheet.Range("E10").Value = '96.39'
$newSheet.Range("F10").Value = '1.34'
$newSheet.Range("G10").Value = '0.0142'
$newSheet.Range("H10").Value = 4

# row 11
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = '515590'
$newSheet.Range("C11").Value = '前海开源中证500等权重ETF'
$newSheet.Range("D11").Value = '0.38'
$newSheet.Range("E11").Value = '95.07'
$newSheet.Range("F11").Value = '0.65'
$newSheet.Range("G11").Value = '0.0025'
$newSheet.Range("H11").Value = 1

# row 12
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = '519620'
$newSheet.Range("C12").Value = '银河君荣灵活配置混合C'
$newSheet.Range("D12").Value = '0.04'
$newSheet.Range("E12").Value = '93.36'
$newSheet.Range("F12").Value = '6.00'
$newSheet.Range("G12").Value = '0.0024'
$newSheet.Range("H12").Value = 5

# row 13
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = '162907'
$newSheet.Range("C13").Value = '泰信中证锐联基本面400指数（LOF）'
$newSheet.Range("D13").Value = '0.23'
$newSheet.Range("E13").Value = '94.61'
$newSheet.Range("F13").Value = '0.85'
$newSheet.Range("G13").Value = '0.0020'
$newSheet.Range("H13").Value = 2

# row 14
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = '519117'
$newSheet.Range("C14").Value = '浦银安盛基本面400指数'
$newSheet.Range("D14").Value = '0.24'
$newSheet.Range("E14").Value = '92.63'
$newSheet.Range("F14").Value = '0.66'
$newSheet.Range("G14").Value = '0.0016'
$newSheet.Range("H14").Value = 6

$newSheet.Range("A1").Select()

# 4. Update "总计": push the existing 2020-Q4 total down to row 3, and insert the
#    new 2022-Q1 total in row 2 (matching the commit's row order).
$oldA = $totalSheet.Range("A2").Value()
$oldB = $totalSheet.Range("B2").Value()
$oldC = $totalSheet.Range("C2").Value()
$oldD = $totalSheet.Range("D2").Value()

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Application.CutCopyMode = $false

$totalSheet.Range("A3").Value = $oldA
$totalSheet.Range("B3").Value = $oldB
$totalSheet.Range("C3").Value = $oldC
$totalSheet.Range("D3").Value = $oldD

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 13
$totalSheet.Range("D2").Value = 0.71

$totalSheet.Range("A1").Select()
